# Add a new worksheet ("Sheet3") at the end of the workbook (after Sheet2)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Header row - write in this order so the shared-string table gets the
# same ordering as the target workbook (Mw drug, Mw co2, S, rho, y1)
$ws3.Range("H1").Value = "Mw drug"
$ws3.Range("I1").Value = "Mw co2"
$ws3.Range("A1").Value = "S"
$ws3.Range("B1").Value = "rho"
$ws3.Range("C1").Value = "y1"

# Constants used by the formulas
$ws3.Range("H2").Value = 598.5
$ws3.Range("I2").Value = 44.01

# Raw data (S, rho) for rows 2-7
$ws3.Range("A2").Value = 0.021
$ws3.Range("B2").Value = 769
$ws3.Range("A3").Value = 0.026
$ws3.Range("B3").Value = 817
$ws3.Range("A4").Value = 0.029
$ws3.Range("B4").Value = 849
$ws3.Range("A5").Value = 0.039
$ws3.Range("B5").Value = 875
$ws3.Range("A6").Value = 0.046
$ws3.Range("B6").Value = 896
$ws3.Range("A7").Value = 0.066
$ws3.Range("B7").Value = 914

# y1 (C) and C*1e6 (D) formulas - row 2 entered on its own first ...
$ws3.Range("C2").Formula = "=(A2*44.01)/((B2*598.5) + (A2*44.01))"
$ws3.Range("D2").Formula = "=C2*1000000"

# ... then rows 3-7 filled together as a single range formula entry
$ws3.Range("C3:C7").Formula = "=(A3*44.01)/((B3*598.5) + (A3*44.01))"
$ws3.Range("D3:D7").Formula = "=C3*1000000"

# Column C width (bestFit, width 12) to match the authored sheet
$ws3.Columns("C").ColumnWidth = 12

# Selection left on the sheet
$null = $ws3.Range("B2:B7").Select()

# Make Sheet3 the active/visible tab
$ws3.Activate()

# Scroll Sheet2 so row 43 is the top visible row (best effort - matches
# the author re-visiting that sheet before saving)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1

# Leave Sheet3 active/selected, as in the saved workbook
$ws3.Activate()
